$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161827087402344
$ws.Range("B1").Value = 2.414998531341553
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.380990028381348
$ws.Range("E1").Value = 1.232399582862854
